# Update town close year columns (#769)
#
# Relabel the "2023 ..." / "2024 ..." year-qualified headers (columns I1:T1)
# to the generic "Prior Year ..." / "Curr. Year ..." labels (same column
# positions, new text only), then apply the PR's reworked column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "Prior Year LMV"
$ws.Range("J1").Value = "Prior Year BMV"
$ws.Range("K1").Value = "Prior Year Total MV"
$ws.Range("L1").Value = "Prior Year LAV"
$ws.Range("M1").Value = "Prior Year BAV"
$ws.Range("N1").Value = "Prior Year Total AV"
$ws.Range("O1").Value = "Curr. Year LMV"
$ws.Range("P1").Value = "Curr. Year BMV"
$ws.Range("Q1").Value = "Curr. Year Total MV"
$ws.Range("R1").Value = "Curr. Year LAV"
$ws.Range("S1").Value = "Curr. Year BAV"
$ws.Range("T1").Value = "Curr. Year Total AV"

# New column widths (declutter / widen columns so headers aren't truncated)
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 8.6667
$ws.Columns.Item(3).ColumnWidth = 13.3333
$ws.Columns.Item(4).ColumnWidth = 7.3333
$ws.Columns.Item(5).ColumnWidth = 8
$ws.Columns.Item(6).ColumnWidth = 23.3333
$ws.Columns.Item(7).ColumnWidth = 19.6667
$ws.Columns.Item(8).ColumnWidth = 16
$ws.Columns.Item(9).ColumnWidth = 15.8333
$ws.Columns.Item(10).ColumnWidth = 16.3333
$ws.Columns.Item(11).ColumnWidth = 20
$ws.Columns.Item(12).ColumnWidth = 16
$ws.Columns.Item(13).ColumnWidth = 15.5
$ws.Columns.Item(14).ColumnWidth = 19.8333
$ws.Columns.Item(15).ColumnWidth = 16
$ws.Columns.Item(16).ColumnWidth = 17.6667
$ws.Columns.Item(17).ColumnWidth = 21.3333
$ws.Columns.Item(18).ColumnWidth = 15.6667
$ws.Columns.Item(19).ColumnWidth = 16.5
$ws.Columns.Item(20).ColumnWidth = 16.5
